$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# The "Previous Doc" column (AX) for every data row (2-15) changes from
# "3471967324" to "9174367677". Stage the new value as real text (not a
# number) on a scratch cell, then copy just the value/type into each
# destination cell so the existing cell style (s="1") is preserved.
$tmp = $ws3.Cells.Item(1, 1)
$tmp.NumberFormat = "@"
$tmp.Value = "9174367677"
$tmp.Copy()
for ($r = 2; $r -le 15; $r++) {
    $ws1.Cells.Item($r, 50).PasteSpecial(-4163)
}
$tmp.Clear()

# Register an extra (unused) font entry in the workbook's style table,
# matching the additional font introduced upstream.
$st = $wb.Styles.Item("Normal 2")
$st.Font.Size = 12
